$wb = $excel.ActiveWorkbook

# The 5 newly-mastered words (same set appended to "Przerobione" and
# refreshed into "5 losowych").
$newWords = @(
    @{ Id = 463;  Hanzi = "经历"; Pinyin = "jīnglì";   English = "experience" },
    @{ Id = 1119; Hanzi = "真";   Pinyin = "zhēn";     English = "really" },
    @{ Id = 889;  Hanzi = "网球"; Pinyin = "wǎngqiú";  English = "Tennis" },
    @{ Id = 195;  Hanzi = "动作"; Pinyin = "dòngzuò";  English = "action" },
    @{ Id = 785;  Hanzi = "手表"; Pinyin = "shǒubiǎo"; English = "Wrist Watch" }
)

# --- Append the 5 words to the bottom of "Przerobione" (rows 47-51) ---
$przerobione = $wb.Worksheets.Item("Przerobione")
$startRow = 47
for ($i = 0; $i -lt $newWords.Count; $i++) {
    $row = $startRow + $i
    $word = $newWords[$i]
    $przerobione.Cells.Item($row, 1).Value = $word.Id
    $przerobione.Cells.Item($row, 2).Value = $word.Hanzi
    $przerobione.Cells.Item($row, 3).Value = $word.Pinyin
    $przerobione.Cells.Item($row, 4).Value = $word.English
}

# --- Refresh "5 losowych" (rows 2-6) with the same 5 words ---
$losowych = $wb.Worksheets.Item("5 losowych")
$startRow2 = 2
for ($i = 0; $i -lt $newWords.Count; $i++) {
    $row = $startRow2 + $i
    $word = $newWords[$i]
    $losowych.Cells.Item($row, 1).Value = $word.Id
    $losowych.Cells.Item($row, 2).Value = $word.Hanzi
    $losowych.Cells.Item($row, 3).Value = $word.Pinyin
    $losowych.Cells.Item($row, 4).Value = $word.English
}
